$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.069761271434878
$ws.Cells.Item(2, 4).Value = 1.067785388946207
$ws.Cells.Item(2, 5).Value = 1.07326956157212
$ws.Cells.Item(2, 6).Value = 1.082819363328973
$ws.Cells.Item(2, 9).Value = 1.02359499962809
$ws.Cells.Item(2, 10).Value = 1.074694357669287
$ws.Cells.Item(2, 11).Value = 1.070492320365493
$ws.Cells.Item(2, 12).Value = 1.075961869338504
$ws.Cells.Item(2, 13).Value = 1.085486588178511
$ws.Cells.Item(2, 14).Value = 1.076220546926343
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.07378490091116
$ws.Cells.Item(3, 4).Value = 1.071509064804203
$ws.Cells.Item(3, 5).Value = 1.076957498081293
$ws.Cells.Item(3, 6).Value = 1.086771039950969
$ws.Cells.Item(3, 9).Value = 1.023504579208684
$ws.Cells.Item(3, 10).Value = 1.078360888699482
$ws.Cells.Item(3, 11).Value = 1.074024634390174
$ws.Cells.Item(3, 12).Value = 1.079459663349571
$ws.Cells.Item(3, 13).Value = 1.089249426160168
$ws.Cells.Item(3, 14).Value = 1.079892284851158
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.076362090783609
$ws.Cells.Item(4, 4).Value = 1.073893282188317
$ws.Cells.Item(4, 5).Value = 1.079318663033241
$ws.Cells.Item(4, 6).Value = 1.089302230544024
$ws.Cells.Item(4, 9).Value = 1.023443335729026
$ws.Cells.Item(4, 10).Value = 1.080707745567631
$ws.Cells.Item(4, 11).Value = 1.076285050242306
$ws.Cells.Item(4, 12).Value = 1.081697799192502
$ws.Cells.Item(4, 13).Value = 1.09165847422222
$ws.Cells.Item(4, 14).Value = 1.082242474525248
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.077439446749066
$ws.Cells.Item(5, 4).Value = 1.074889764089804
$ws.Cells.Item(5, 5).Value = 1.080305469383531
$ws.Cells.Item(5, 6).Value = 1.090360376086116
$ws.Cells.Item(5, 9).Value = 1.023416932628352
$ws.Cells.Item(5, 10).Value = 1.081688428980044
$ws.Cells.Item(5, 11).Value = 1.077229484796218
$ws.Cells.Item(5, 12).Value = 1.08263287913079
$ws.Cells.Item(5, 13).Value = 1.092665275982175
$ws.Cells.Item(5, 14).Value = 1.083224550620593
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.077619988369254
$ws.Cells.Item(6, 4).Value = 1.075056740895716
$ws.Cells.Item(6, 5).Value = 1.080470822488941
$ws.Cells.Item(6, 6).Value = 1.090537699503253
$ws.Cells.Item(6, 9).Value = 1.023412460900047
$ws.Cells.Item(6, 10).Value = 1.081852747741767
$ws.Cells.Item(6, 11).Value = 1.077387722367406
$ws.Cells.Item(6, 12).Value = 1.082789546670653
$ws.Cells.Item(6, 13).Value = 1.092833978650151
$ws.Cells.Item(6, 14).Value = 1.083389102733805
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.076376510130172
$ws.Cells.Item(7, 4).Value = 1.073906619924006
$ws.Cells.Item(7, 5).Value = 1.079331871424509
$ws.Cells.Item(7, 6).Value = 1.089316392706307
$ws.Cells.Item(7, 9).Value = 1.02344298551034
$ws.Cells.Item(7, 10).Value = 1.080720872565408
$ws.Cells.Item(7, 11).Value = 1.076297692531983
$ws.Cells.Item(7, 12).Value = 1.081710316440782
$ws.Cells.Item(7, 13).Value = 1.091671950315747
$ws.Cells.Item(7, 14).Value = 1.082255620164868
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.071126668725622
$ws.Cells.Item(8, 4).Value = 1.06904917432753
$ws.Cells.Item(8, 5).Value = 1.074521252353219
$ws.Cells.Item(8, 6).Value = 1.084160324290593
$ws.Cells.Item(8, 9).Value = 1.023565007130014
$ws.Cells.Item(8, 10).Value = 1.075938910077685
$ws.Cells.Item(8, 11).Value = 1.071691425488774
$ws.Cells.Item(8, 12).Value = 1.077149295479841
$ws.Cells.Item(8, 13).Value = 1.086763716044522
$ws.Cells.Item(8, 14).Value = 1.077466866741913
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.061663718938178
$ws.Cells.Item(9, 4).Value = 1.060286977572865
$ws.Cells.Item(9, 5).Value = 1.06584221789996
$ws.Cells.Item(9, 6).Value = 1.074867182822559
$ws.Cells.Item(9, 9).Value = 1.023759139049229
$ws.Cells.Item(9, 10).Value = 1.067306791918144
$ws.Cells.Item(9, 11).Value = 1.063372311689228
$ws.Cells.Item(9, 12).Value = 1.068910423956386
$ws.Cells.Item(9, 13).Value = 1.077907974389667
$ws.Cells.Item(9, 14).Value = 1.068822489984467
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.055197854270927
$ws.Cells.Item(10, 4).Value = 1.054295611293696
$ws.Cells.Item(10, 5).Value = 1.059906854158194
$ws.Cells.Item(10, 6).Value = 1.068518011982952
$ws.Cells.Item(10, 9).Value = 1.023874603872166
$ws.Cells.Item(10, 10).Value = 1.061400184710702
$ws.Cells.Item(10, 11).Value = 1.057677090588836
$ws.Cells.Item(10, 12).Value = 1.063269155963708
$ws.Cells.Item(10, 13).Value = 1.071851313715718
$ws.Cells.Item(10, 14).Value = 1.062907494717295
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.052357348190432
$ws.Cells.Item(11, 4).Value = 1.051662558193425
$ws.Cells.Item(11, 5).Value = 1.057298213136342
$ws.Cells.Item(11, 6).Value = 1.065728960723682
$ws.Cells.Item(11, 9).Value = 1.023921309042974
$ws.Cells.Item(11, 10).Value = 1.058803370737039
$ws.Cells.Item(11, 11).Value = 1.055172552928917
$ws.Cells.Item(11, 12).Value = 1.060788115781959
$ws.Cells.Item(11, 13).Value = 1.069189243321588
$ws.Cells.Item(11, 14).Value = 1.06030699296993
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.05129583071583
$ws.Cells.Item(12, 4).Value = 1.050678419744675
$ws.Cells.Item(12, 5).Value = 1.05632316946924
$ws.Cells.Item(12, 6).Value = 1.064686703683364
$ws.Cells.Item(12, 9).Value = 1.023938164268118
$ws.Cells.Item(12, 10).Value = 1.057832623876867
$ws.Cells.Item(12, 11).Value = 1.054236202368447
$ws.Cells.Item(12, 12).Value = 1.059860515943885
$ws.Cells.Item(12, 13).Value = 1.06819421034859
$ws.Cells.Item(12, 14).Value = 1.05933486753787
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.051523825717757
$ws.Cells.Item(13, 4).Value = 1.050889801750513
$ws.Cells.Item(13, 5).Value = 1.056532599346114
$ws.Cells.Item(13, 6).Value = 1.064910560422093
$ws.Cells.Item(13, 9).Value = 1.023934571058996
$ws.Cells.Item(13, 10).Value = 1.058041136471595
$ws.Cells.Item(13, 11).Value = 1.054437331228016
$ws.Cells.Item(13, 12).Value = 1.060059766676768
$ws.Cells.Item(13, 13).Value = 1.068407934598042
$ws.Cells.Item(13, 14).Value = 1.059543676244399
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.052269735750752
$ws.Cells.Item(14, 4).Value = 1.051581335230925
$ws.Cells.Item(14, 5).Value = 1.057217741383971
$ws.Cells.Item(14, 6).Value = 1.065642937289826
$ws.Cells.Item(14, 9).Value = 1.023922712353274
$ws.Cells.Item(14, 10).Value = 1.058723256138755
$ws.Cells.Item(14, 11).Value = 1.055095279023094
$ws.Cells.Item(14, 12).Value = 1.060711564736175
$ws.Cells.Item(14, 13).Value = 1.069107122210945
$ws.Cells.Item(14, 14).Value = 1.060226764599727
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.052728454518375
$ws.Cells.Item(15, 4).Value = 1.052006594134751
$ws.Cells.Item(15, 5).Value = 1.057639065973708
$ws.Cells.Item(15, 6).Value = 1.066093337684767
$ws.Cells.Item(15, 9).Value = 1.023915340507297
$ws.Cells.Item(15, 10).Value = 1.05914270571262
$ws.Cells.Item(15, 11).Value = 1.055499851784162
$ws.Cells.Item(15, 12).Value = 1.061112351524428
$ws.Cells.Item(15, 13).Value = 1.069537081545364
$ws.Cells.Item(15, 14).Value = 1.060646809840098
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.055385485092558
$ws.Cells.Item(16, 4).Value = 1.054469518008761
$ws.Cells.Item(16, 5).Value = 1.060079144370744
$ws.Cells.Item(16, 6).Value = 1.068702248157306
$ws.Cells.Item(16, 9).Value = 1.023871434968932
$ws.Cells.Item(16, 10).Value = 1.06157167663451
$ws.Cells.Item(16, 11).Value = 1.057842474877153
$ws.Cells.Item(16, 12).Value = 1.063432983854572
$ws.Cells.Item(16, 13).Value = 1.072027130160847
$ws.Cells.Item(16, 14).Value = 1.063079230179304
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.057041050305959
$ws.Cells.Item(17, 4).Value = 1.056003874084446
$ws.Cells.Item(17, 5).Value = 1.061599215428059
$ws.Cells.Item(17, 6).Value = 1.070327882768646
$ws.Cells.Item(17, 9).Value = 1.023843013862359
$ws.Cells.Item(17, 10).Value = 1.063084611336228
$ws.Cells.Item(17, 11).Value = 1.059301451397465
$ws.Cells.Item(17, 12).Value = 1.064878204044377
$ws.Cells.Item(17, 13).Value = 1.073578298551096
$ws.Cells.Item(17, 14).Value = 1.064594313421835
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.058002803626237
$ws.Cells.Item(18, 4).Value = 1.056895118798761
$ws.Cells.Item(18, 5).Value = 1.062482143439312
$ws.Cells.Item(18, 6).Value = 1.071272266771262
$ws.Cells.Item(18, 9).Value = 1.023826118601775
$ws.Cells.Item(18, 10).Value = 1.063963316726034
$ws.Cells.Item(18, 11).Value = 1.060148755128853
$ws.Cells.Item(18, 12).Value = 1.065717496842228
$ws.Cells.Item(18, 13).Value = 1.074479278086792
$ws.Cells.Item(18, 14).Value = 1.065474266674083
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.058330082532017
$ws.Cells.Item(19, 4).Value = 1.057198387761433
$ws.Cells.Item(19, 5).Value = 1.06278257917382
$ws.Cells.Item(19, 6).Value = 1.07159363797398
$ws.Cells.Item(19, 9).Value = 1.023820303822167
$ws.Cells.Item(19, 10).Value = 1.064262302532385
$ws.Cells.Item(19, 11).Value = 1.060437045600963
$ws.Cells.Item(19, 12).Value = 1.066003057911241
$ws.Cells.Item(19, 13).Value = 1.074785854370916
$ws.Cells.Item(19, 14).Value = 1.065773677074573
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.0568638304078
$ws.Cells.Item(20, 4).Value = 1.055839638902447
$ws.Cells.Item(20, 5).Value = 1.061436511270899
$ws.Cells.Item(20, 6).Value = 1.07015386490463
$ws.Cells.Item(20, 9).Value = 1.023846096030143
$ws.Cells.Item(20, 10).Value = 1.062922679076665
$ws.Cells.Item(20, 11).Value = 1.059145300908942
$ws.Cells.Item(20, 12).Value = 1.064723528128733
$ws.Cells.Item(20, 13).Value = 1.073412266995193
$ws.Cells.Item(20, 14).Value = 1.064432151199889
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.052050263780348
$ws.Cells.Item(21, 4).Value = 1.051377866732292
$ws.Cells.Item(21, 5).Value = 1.057016154237538
$ws.Cells.Item(21, 6).Value = 1.065427446346162
$ws.Cells.Item(21, 9).Value = 1.02392621804881
$ws.Cells.Item(21, 10).Value = 1.058522561683257
$ws.Cells.Item(21, 11).Value = 1.0549016991743
$ws.Cells.Item(21, 12).Value = 1.060519795184712
$ws.Cells.Item(21, 13).Value = 1.068901403005153
$ws.Cells.Item(21, 14).Value = 1.060025785135082
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.048986436919783
$ws.Cells.Item(22, 4).Value = 1.048537100465997
$ws.Cells.Item(22, 5).Value = 1.054201585194992
$ws.Cells.Item(22, 6).Value = 1.062419272844665
$ws.Cells.Item(22, 9).Value = 1.023973741382756
$ws.Cells.Item(22, 10).Value = 1.055720159730901
$ws.Cells.Item(22, 11).Value = 1.052198408252213
$ws.Cells.Item(22, 12).Value = 1.057841704016115
$ws.Cells.Item(22, 13).Value = 1.066029096097794
$ws.Cells.Item(22, 14).Value = 1.057219403450514
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.050614275254111
$ws.Cells.Item(23, 4).Value = 1.050046504798131
$ws.Cells.Item(23, 5).Value = 1.055697086010733
$ws.Cells.Item(23, 6).Value = 1.064017523733382
$ws.Cells.Item(23, 9).Value = 1.023948818288664
$ws.Cells.Item(23, 10).Value = 1.057209264314164
$ws.Cells.Item(23, 11).Value = 1.053634902454559
$ws.Cells.Item(23, 12).Value = 1.059264825958132
$ws.Cells.Item(23, 13).Value = 1.067555286184473
$ws.Cells.Item(23, 14).Value = 1.058710622733087
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.056943920512945
$ws.Cells.Item(24, 4).Value = 1.055913861185716
$ws.Cells.Item(24, 5).Value = 1.061510041702207
$ws.Cells.Item(24, 6).Value = 1.070232507874908
$ws.Cells.Item(24, 9).Value = 1.023844704313072
$ws.Cells.Item(24, 10).Value = 1.062995860906222
$ws.Cells.Item(24, 11).Value = 1.059215869985504
$ws.Cells.Item(24, 12).Value = 1.064793430869575
$ws.Cells.Item(24, 13).Value = 1.073487301203203
$ws.Cells.Item(24, 14).Value = 1.064505436956038
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.064136719599562
$ws.Cells.Item(25, 4).Value = 1.062577606536118
$ws.Cells.Item(25, 5).Value = 1.068111258050262
$ws.Cells.Item(25, 6).Value = 1.077295700824986
$ws.Cells.Item(25, 9).Value = 1.02371141742819
$ws.Cells.Item(25, 10).Value = 1.069564138732758
$ws.Cells.Item(25, 11).Value = 1.065548292893455
$ws.Cells.Item(25, 12).Value = 1.071065588445654
$ws.Cells.Item(25, 13).Value = 1.080223287716545
$ws.Cells.Item(25, 14).Value = 1.071083042490479
